$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 290
$ws1.Range("F6").Value = 112
$ws1.Range("F7").Value = 290
$ws1.Range("F9").Value = 2038
$ws1.Range("F10").Value = 358
$ws1.Range("F11").Value = 4855
$ws1.Range("F12").Value = 94
$ws1.Range("F13").Value = 339

# Sheet "全部类型" - update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 290
$ws4.Range("F8").Value = 112
$ws4.Range("F9").Value = 290
$ws4.Range("F13").Value = 2038
$ws4.Range("F14").Value = 358
$ws4.Range("F15").Value = 4855
$ws4.Range("F16").Value = 94
$ws4.Range("F17").Value = 339
